# New operations and battles: mark the last 5 entries (Battle of Okinawa,
# Battle for the Philippine Islands / Battle of the Philippine sea,
# Battle of Leyte Gulf, Battle of Iwo Jima, Operation Market Garden)
# as done ("x") in column D, and move the active selection to D30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D30").Value = "x"
$ws.Range("D31").Value = "x"
$ws.Range("D32").Value = "x"
$ws.Range("D33").Value = "x"
$ws.Range("D34").Value = "x"

$ws.Range("D30").Select()
